# Femacal de La Calera - Chirimoya: add a new week's price report (2023-10-12)
# at the top of the Chirimoya data block (rows 342-425), pushing the
# existing rows down by 3 (to rows 345-428).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 342. Excel automatically shifts the
# existing rows 342:425 down to 345:428, preserving their values/styles.
$ws.Rows("342:344").Insert()

# New date for this week's report: 2023-10-12
$newDate = 45211

# Row 342: "Especial" quality
$ws.Cells.Item(342, 1).Value = 3
$ws.Cells.Item(342, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(342, 3).Value = "Coquimbo"
$ws.Cells.Item(342, 4).Value = $newDate
$ws.Cells.Item(342, 5).Value = 5
$ws.Cells.Item(342, 6).Value = "Fruta"
$ws.Cells.Item(342, 7).Value = 100107
$ws.Cells.Item(342, 8).Value = "Otros"
$ws.Cells.Item(342, 9).Value = 100107002
$ws.Cells.Item(342, 10).Value = "Chirimoya"
$ws.Cells.Item(342, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(342, 12).Value = "Especial"
$ws.Cells.Item(342, 13).Value = 101
$ws.Cells.Item(342, 14).Value = 28000
$ws.Cells.Item(342, 15).Value = 30000
$ws.Cells.Item(342, 16).Value = 29069
$ws.Cells.Item(342, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(342, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(342, 19).Value = 2907
$ws.Cells.Item(342, 20).Value = 10

# Row 343: "Primera" quality
$ws.Cells.Item(343, 1).Value = 3
$ws.Cells.Item(343, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(343, 3).Value = "Coquimbo"
$ws.Cells.Item(343, 4).Value = $newDate
$ws.Cells.Item(343, 5).Value = 5
$ws.Cells.Item(343, 6).Value = "Fruta"
$ws.Cells.Item(343, 7).Value = 100107
$ws.Cells.Item(343, 8).Value = "Otros"
$ws.Cells.Item(343, 9).Value = 100107002
$ws.Cells.Item(343, 10).Value = "Chirimoya"
$ws.Cells.Item(343, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(343, 12).Value = "Primera"
$ws.Cells.Item(343, 13).Value = 97
$ws.Cells.Item(343, 14).Value = 25000
$ws.Cells.Item(343, 15).Value = 26000
$ws.Cells.Item(343, 16).Value = 25412
$ws.Cells.Item(343, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(343, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(343, 19).Value = 2541
$ws.Cells.Item(343, 20).Value = 10

# Row 344: "Segunda" quality
$ws.Cells.Item(344, 1).Value = 3
$ws.Cells.Item(344, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(344, 3).Value = "Coquimbo"
$ws.Cells.Item(344, 4).Value = $newDate
$ws.Cells.Item(344, 5).Value = 5
$ws.Cells.Item(344, 6).Value = "Fruta"
$ws.Cells.Item(344, 7).Value = 100107
$ws.Cells.Item(344, 8).Value = "Otros"
$ws.Cells.Item(344, 9).Value = 100107002
$ws.Cells.Item(344, 10).Value = "Chirimoya"
$ws.Cells.Item(344, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(344, 12).Value = "Segunda"
$ws.Cells.Item(344, 13).Value = 68
$ws.Cells.Item(344, 14).Value = 22000
$ws.Cells.Item(344, 15).Value = 22000
$ws.Cells.Item(344, 16).Value = 22000
$ws.Cells.Item(344, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(344, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(344, 19).Value = 2200
$ws.Cells.Item(344, 20).Value = 10
